$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 3200.25
$ws.Range("I2").Value = 324.83334
$ws.Range("K2").Value = 324.83334
$ws.Range("M2").Value = -211.83334

# Row 4
$ws.Range("H4").Value = 15150
$ws.Range("J4").Value = 15150
$ws.Range("L4").Value = 15150
$ws.Range("N4").Value = -15378

# Row 18
$ws.Range("H18").Value = 1110.238
$ws.Range("I18").Value = 1137.1
$ws.Range("K18").Value = 1137.1
$ws.Range("M18").Value = -853.0999999999999

# Row 33
$ws.Range("H33").Value = 143.13333
$ws.Range("I33").Value = 142.46153
$ws.Range("K33").Value = 142.46153
$ws.Range("M33").Value = 86.53846999999999

# Row 94
$ws.Range("H94").Value = 931
$ws.Range("I94").Value = 896.5
$ws.Range("K94").Value = 896.5
$ws.Range("M94").Value = -445.5

# Row 135
$ws.Range("H135").Value = 4509.231
$ws.Range("I135").Value = 3715.3333
$ws.Range("K135").Value = 33437.9997
$ws.Range("M135").Value = -30902.9997

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 4805.2856
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4805.2856
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4805.2856
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5443.2856

# Row 88
$ws.Range("H88").Value = 2422.1428
$ws.Range("J88").Value = 2241.5
$ws.Range("L88").Value = 2241.5
$ws.Range("N88").Value = -3053.5

# Row 91
$ws.Range("H91").Value = 2422.1428
$ws.Range("J91").Value = 2241.5
$ws.Range("L91").Value = 2241.5
$ws.Range("N91").Value = -5049.5

# Row 102
$ws.Range("H102").Value = 799
$ws.Range("I102").Value = 799
$ws.Range("K102").Value = 799
$ws.Range("M102").Value = 823

# Row 110
$ws.Range("H110").Value = 1122.1111
$ws.Range("I110").Value = 1063.1765
$ws.Range("K110").Value = 1063.1765
$ws.Range("M110").Value = 981.8235

# Row 132
$ws.Range("H132").Value = 1994.5
$ws.Range("I132").Value = 1994
$ws.Range("K132").Value = 5982
$ws.Range("M132").Value = -3452

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 260.5
$ws.Range("I22").Value = 272.2
$ws.Range("K22").Value = 272.2
$ws.Range("M22").Value = -99.19999999999999

# Row 94
$ws.Range("H94").Value = 2999
$ws.Range("I94").Value = 2999
$ws.Range("K94").Value = 2999
$ws.Range("M94").Value = -2548

# Row 105
$ws.Range("H105").Value = 2799.3333
$ws.Range("I105").Value = 2799.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2799.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1052.3333
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 2999.077
$ws.Range("I105").Value = 1998.2858
$ws.Range("J105").Value = 4166.6665
$ws.Range("K105").Value = 1998.2858
$ws.Range("L105").Value = 4166.6665
$ws.Range("M105").Value = -251.2858000000001
$ws.Range("N105").Value = -7660.6665

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1316.3334
$ws.Range("I113").Value = 650
$ws.Range("K113").Value = 1950
$ws.Range("M113").Value = 220

# Row 122
$ws.Range("H122").Value = 887.6667
$ws.Range("J122").Value = 887.6667
$ws.Range("L122").Value = 7989.0003
$ws.Range("N122").Value = -12889.0003

# Row 129
$ws.Range("H129").Value = 1668505.5
$ws.Range("I129").Value = 1333.3334
$ws.Range("J129").Value = 3335677.8
$ws.Range("K129").Value = 4000.0002
$ws.Range("L129").Value = 10007033.4
$ws.Range("M129").Value = 999.9998000000001
$ws.Range("N129").Value = -10017033.4

# Row 131
$ws.Range("H131").Value = 668557.2
$ws.Range("J131").Value = 1113382.6
$ws.Range("L131").Value = 3340147.8
$ws.Range("N131").Value = -3350227.8

# Row 134
$ws.Range("H134").Value = 1137.4
$ws.Range("I134").Value = 1072
$ws.Range("J134").Value = 1399
$ws.Range("K134").Value = 3216
$ws.Range("L134").Value = 4197
$ws.Range("M134").Value = 1854
$ws.Range("N134").Value = -14337

# Row 138
$ws.Range("H138").Value = 1976.8
$ws.Range("I138").Value = 1871.25
$ws.Range("K138").Value = 5613.75
$ws.Range("M138").Value = -473.75

# Row 140
$ws.Range("H140").Value = 4716.5
$ws.Range("I140").Value = 309.7143
$ws.Range("J140").Value = 14999
$ws.Range("K140").Value = 929.1428999999999
$ws.Range("L140").Value = 44997
$ws.Range("M140").Value = 4250.8571
$ws.Range("N140").Value = -55357

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 47500
$ws.Range("J40").Value = 47500
$ws.Range("L40").Value = 47500
$ws.Range("N40").Value = -47802

# Row 122
$ws.Range("H122").Value = 3898.8572
$ws.Range("I122").Value = 3898.8572
$ws.Range("K122").Value = 11696.5716
$ws.Range("M122").Value = -9246.571599999999

# Row 126
$ws.Range("H126").Value = 8487.5
$ws.Range("J126").Value = 9000
$ws.Range("L126").Value = 27000
$ws.Range("N126").Value = -31940

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7956.067
$ws.Range("I7").Value = 4889.5
$ws.Range("K7").Value = 4889.5
$ws.Range("M7").Value = -4777.5

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 46
$ws.Range("H46").Value = 1469.6154
$ws.Range("I46").Value = 1315.2222
$ws.Range("J46").Value = 1551.3529
$ws.Range("K46").Value = 1315.2222
$ws.Range("L46").Value = 1551.3529
$ws.Range("M46").Value = -1127.2222
$ws.Range("N46").Value = -1927.3529

# Row 82
$ws.Range("H82").Value = 2932.6667
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2932.6667
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2932.6667
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3654.6667

# Row 85
$ws.Range("H85").Value = 2932.6667
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2932.6667
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2932.6667
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5428.6667

# Row 100
$ws.Range("H100").Value = 3681.1052
$ws.Range("I100").Value = 3677.6875
$ws.Range("J100").Value = 3699.3333
$ws.Range("K100").Value = 3677.6875
$ws.Range("L100").Value = 3699.3333
$ws.Range("M100").Value = -3136.6875
$ws.Range("N100").Value = -4781.3333

# Row 122
$ws.Range("H122").Value = 6309.607
$ws.Range("I122").Value = 4799
$ws.Range("K122").Value = 14397
$ws.Range("M122").Value = -11947

# Row 126
$ws.Range("H126").Value = 7956.067
$ws.Range("I126").Value = 4889.5
$ws.Range("K126").Value = 14668.5
$ws.Range("M126").Value = -12198.5

# Row 133
$ws.Range("H133").Value = 90296
$ws.Range("I133").Value = 90296
$ws.Range("K133").Value = 90296
$ws.Range("M133").Value = -87766

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 9284.125
$ws.Range("I122").Value = 8644
$ws.Range("K122").Value = 25932
$ws.Range("M122").Value = -23482

# Row 126
$ws.Range("H126").Value = 2572.5557
$ws.Range("I126").Value = 1550.8
$ws.Range("K126").Value = 4652.4
$ws.Range("M126").Value = -2182.4

# Row 132
$ws.Range("H132").Value = 2815.3333
$ws.Range("I132").Value = 2815.3333
$ws.Range("K132").Value = 8445.999899999999
$ws.Range("M132").Value = -5915.999899999999
